# Update cryptocurrency price/volume figures per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$priceStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.377.10"
$ws.Range("D2").Style = $priceStyle
$ws.Range("E2").Value = "  -1.02%  "

# Row 3
$priceStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.884.03"
$ws.Range("D3").Style = $priceStyle
$ws.Range("E3").Value = "  -1.43%  "

# Row 4
$priceStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = $priceStyle
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$priceStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.36"
$ws.Range("D5").Style = $priceStyle
$ws.Range("E5").Value = "  -0.51%  "

# Row 6
$priceStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = $priceStyle
$ws.Range("E6").Value = "  -0.20%  "

# Row 7
$priceStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4830"
$ws.Range("D7").Style = $priceStyle
$ws.Range("E7").Value = "  -1.97%  "

# Row 8
$priceStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2883"
$ws.Range("D8").Style = $priceStyle
$ws.Range("E8").Value = "  -2.89%  "

# Row 9
$priceStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06597"
$ws.Range("D9").Style = $priceStyle
$ws.Range("E9").Value = "  -2.58%  "

# Row 10
$priceStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.897.23"
$ws.Range("D10").Style = $priceStyle
$ws.Range("E10").Value = "  -1.05%  "

# Row 11
$priceStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.93"
$ws.Range("D11").Style = $priceStyle
$ws.Range("E11").Value = "  -1.29%  "

# Row 12
$priceStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07397"
$ws.Range("D12").Style = $priceStyle
$ws.Range("E12").Value = "  +0.50%  "

# Row 13
$priceStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.162"
$ws.Range("D13").Style = $priceStyle
$ws.Range("E13").Value = "  -0.12%  "

# Row 14
$priceStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.56"
$ws.Range("D14").Style = $priceStyle
$ws.Range("E14").Value = "  +0.03%  "

# Row 15
$priceStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6606"
$ws.Range("D15").Style = $priceStyle
$ws.Range("E15").Value = "  -1.78%  "

# Row 16
$priceStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.326.94"
$ws.Range("D16").Style = $priceStyle
$ws.Range("E16").Value = "  -1.05%  "

# Row 17
$ws.Range("E17").Value = "  -0.16%  "

# Row 18
$priceStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007753"
$ws.Range("D18").Style = $priceStyle
$ws.Range("E18").Value = "  -2.41%  "

# Row 19
$priceStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9996"
$ws.Range("D19").Style = $priceStyle
$ws.Range("E19").Value = "  -0.14%  "

# Row 20
$priceStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.436"
$ws.Range("D20").Style = $priceStyle
$ws.Range("E20").Value = "  +0.94%  "

# Row 21
$priceStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.126.47"
$ws.Range("D21").Style = $priceStyle
$ws.Range("E21").Value = "  -1.08%  "

# Row 22
$priceStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9989"
$ws.Range("D22").Style = $priceStyle
$ws.Range("E22").Value = "  -0.40%  "

# Row 23
$priceStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "197.56"
$ws.Range("D23").Style = $priceStyle
$ws.Range("E23").Value = "  -0.86%  "

# Row 24
$priceStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.164"
$ws.Range("D24").Style = $priceStyle
$ws.Range("E24").Value = "  -1.76%  "

# Row 25
$priceStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.416"
$ws.Range("D25").Style = $priceStyle
$ws.Range("E25").Value = "  -2.56%  "

# Row 26
$priceStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.15"
$ws.Range("D26").Style = $priceStyle
$ws.Range("E26").Value = "  +0.13%  "

# Row 27
$priceStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.21"
$ws.Range("D27").Style = $priceStyle
$ws.Range("E27").Value = "  -2.50%  "

# Row 28
$priceStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.939"
$ws.Range("D28").Style = $priceStyle
$ws.Range("E28").Value = "  -0.68%  "

# Row 29
$priceStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.434"
$ws.Range("D29").Style = $priceStyle
$ws.Range("E29").Value = "  -3.99%  "

# Row 30
$priceStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.275"
$ws.Range("D30").Style = $priceStyle
$ws.Range("E30").Value = "  -1.79%  "

# Row 31
$priceStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09149"
$ws.Range("D31").Style = $priceStyle
$ws.Range("E31").Value = "  +0.12%  "

# Row 32
$priceStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.055"
$ws.Range("D32").Style = $priceStyle
$ws.Range("E32").Value = "  +0.06%  "

# Row 33
$priceStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05054"
$ws.Range("D33").Style = $priceStyle
$ws.Range("E33").Value = "  -4.65%  "

# Row 34
$priceStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7401"
$ws.Range("D34").Style = $priceStyle
$ws.Range("E34").Value = "  -0.06%  "

# Row 35
$priceStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.151"
$ws.Range("D35").Style = $priceStyle
$ws.Range("E35").Value = "  +3.52%  "

# Row 36
$priceStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.705"
$ws.Range("D36").Style = $priceStyle
$ws.Range("E36").Value = "  -0.55%  "

# Row 37
$priceStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01873"
$ws.Range("D37").Style = $priceStyle
$ws.Range("E37").Value = "  +2.63%  "

# Row 38
$priceStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.634"
$ws.Range("D38").Style = $priceStyle
$ws.Range("E38").Value = "  -3.08%  "

# Row 39
$priceStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9158"
$ws.Range("D39").Style = $priceStyle
$ws.Range("E39").Value = "  -0.69%  "

# Row 40
$priceStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.077"
$ws.Range("D40").Style = $priceStyle
$ws.Range("E40").Value = "  -0.70%  "

# Row 41
$priceStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.974"
$ws.Range("D41").Style = $priceStyle
$ws.Range("E41").Value = "  +0.48%  "

# Row 42
$priceStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.96"
$ws.Range("D42").Style = $priceStyle
$ws.Range("E42").Value = "  -0.01%  "

# Row 43
$priceStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4331"
$ws.Range("D43").Style = $priceStyle
$ws.Range("E43").Value = "  -2.83%  "

# Row 44
$priceStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = $priceStyle
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$priceStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.600"
$ws.Range("D45").Style = $priceStyle
$ws.Range("E45").Value = "  -0.46%  "

# Row 46
$priceStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1347"
$ws.Range("D46").Style = $priceStyle
$ws.Range("E46").Value = "  -3.10%  "

# Row 47
$priceStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.593"
$ws.Range("D47").Style = $priceStyle
$ws.Range("E47").Value = "  +10.58%  "

# Row 48
$priceStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.10"
$ws.Range("D48").Style = $priceStyle
$ws.Range("E48").Value = "  -13.07%  "

# Row 49
$priceStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.881"
$ws.Range("D49").Style = $priceStyle
$ws.Range("E49").Value = "  -3.00%  "

# Row 50
$priceStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.28"
$ws.Range("D50").Style = $priceStyle
$ws.Range("E50").Value = "  -3.33%  "

# Row 51
$priceStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05690"
$ws.Range("D51").Style = $priceStyle
$ws.Range("E51").Value = "  -3.21%  "
